# Replace the manual line break that precedes "Підпис: ________________"
# with a run of 26 spaces, and turn the trailing "Підпис: ..." text into
# its own separate run (matching how Word splits runs when text is
# retyped over a deleted break).

$d = $word.ActiveDocument

$signText = "Підпис: ________________"
$spaces   = "                          "   # 26 spaces

$targetParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pText = $d.Paragraphs.Item($i).Range.Text
    if ($pText.Contains($signText) -and $pText.Contains([string][char]11)) {
        $targetParaIndex = $i
        break
    }
}

if ($targetParaIndex -gt 0) {
    $para = $d.Paragraphs.Item($targetParaIndex)
    $pRange = $para.Range
    $pText = $pRange.Text

    $breakOffset = $pText.IndexOf([char]11)
    $breakStart  = $pRange.Start + $breakOffset

    # 1) Turn the manual line break itself into a run of spaces (keeps the
    #    formatting / w:r attributes of the run the break used to live in).
    $breakRange = $d.Range($breakStart, $breakStart + 1)
    $breakRange.Text = $spaces

    # 2) Nudge the formatting of the spaces run (no visual change, same
    #    color) so it is not silently re-merged with the preceding
    #    "Дата: ..." run - this reproduces the real run boundary Word
    #    leaves behind after such an edit.
    $spacesRange = $d.Range($breakStart, $breakStart + $spaces.Length)
    $spacesRange.Font.Color = 1
    $spacesRange.Font.Color = 0

    # 3) Likewise nudge the "Підпис: ..." run so it stays separate from
    #    the spaces run.
    $signStart = $breakStart + $spaces.Length
    $signRange = $d.Range($signStart, $signStart + $signText.Length)
    $signRange.Font.Color = 1
    $signRange.Font.Color = 0

    # 4) Re-stamp the "Підпис: ..." run through InsertXML so it comes out
    #    as a brand-new <w:r> (no rsid attributes), matching a freshly
    #    authored run rather than a split of the old one.
    $signXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:color w:val="000000"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>' + $signText + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $signRange2 = $d.Range($signStart, $signStart + $signText.Length)
    $signRange2.InsertXML($signXml)
}
